$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (shortened) "Sample ID" query text for cell B3 -- the Tumor and
# Analyte Type columns were removed from the SELECT list.
$newSampleQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001437' AND gi.platform = 'Illumina'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

# Update the query text stored in B3 (this is the only textual edit; the
# shared-strings table reindexes itself as a side effect, same as Excel
# would do: the old string is dropped and the new one appended while the
# untouched "File Name" query slides down one slot).
$ws.Range("B3").Value = $newSampleQuery

# The TsvExcel / WebExcel helper columns (D and E) are no longer needed
# for the Samples and Files rows, so remove those cells entirely.
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# Reflect the author's final cursor position: cell B3 selected, with the
# view scrolled so row 3 is at the top.
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 3
